# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2 through 220) from serial date 46060 to 46061 (i.e. +1 day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C220").Value = 46061
